$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.07794266666666667
$ws.Range("H2").Value = 0.233828
$ws.Range("I2").Value = 0.002827880818927331
$ws.Range("J2").Value = 0.00282788081892733
$ws.Range("M2").Value = 0.06166766666666667
$ws.Range("N2").Value = 0.185003
$ws.Range("O2").Value = 0.3189772891852935
$ws.Range("P2").Value = 0.3189772891852935
$ws.Range("Q2").Value = 0.004806542387111112
$ws.Range("R2").Value = 0.043258881484
$ws.Range("S2").Value = 0.0009020297577605279
$ws.Range("T2").Value = 0.0009020297577605275
# Row 3
$ws.Range("G3").Value = 0.07794266666666667
$ws.Range("H3").Value = 0.233828
$ws.Range("I3").Value = 0.002827880818927331
$ws.Range("J3").Value = 0.00282788081892733
$ws.Range("O3").Value = 0.4045463009579509
$ws.Range("P3").Value = 0.4045463009579509
$ws.Range("Q3").Value = 0.006095947921777778
$ws.Range("R3").Value = 0.054863531296
$ws.Range("S3").Value = 0.001144008724846993
$ws.Range("T3").Value = 0.001144008724846992
# Row 4
$ws.Range("G4").Value = 0.07794266666666667
$ws.Range("H4").Value = 0.233828
$ws.Range("I4").Value = 0.002827880818927331
$ws.Range("J4").Value = 0.00282788081892733
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.05345100000000001
$ws.Range("N4").Value = 0.160353
$ws.Range("O4").Value = 0.2764764098567557
$ws.Range("P4").Value = 0.2764764098567557
$ws.Range("Q4").Value = 0.004166113476000001
$ws.Range("R4").Value = 0.03749502128400001
$ws.Range("S4").Value = 0.0007818423363198107
$ws.Range("T4").Value = 0.0007818423363198105
# Row 5
$ws.Range("I5").Value = 0.9151728997907317
$ws.Range("J5").Value = 0.9151728997907316
$ws.Range("M5").Value = 0.06166766666666667
$ws.Range("N5").Value = 0.185003
$ws.Range("O5").Value = 0.3189772891852935
$ws.Range("P5").Value = 0.3189772891852935
$ws.Range("Q5").Value = 1.555517228639111
$ws.Range("R5").Value = 13.999655057752
$ws.Range("S5").Value = 0.2919193707110919
$ws.Range("T5").Value = 0.2919193707110918
# Row 6
$ws.Range("I6").Value = 0.9151728997907317
$ws.Range("J6").Value = 0.9151728997907316
$ws.Range("O6").Value = 0.4045463009579509
$ws.Range("P6").Value = 0.4045463009579509
$ws.Range("S6").Value = 0.3702298113473019
$ws.Range("T6").Value = 0.3702298113473019
# Row 7
$ws.Range("I7").Value = 0.9151728997907317
$ws.Range("J7").Value = 0.9151728997907316
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.05345100000000001
$ws.Range("N7").Value = 0.160353
$ws.Range("O7").Value = 0.2764764098567557
$ws.Range("P7").Value = 0.2764764098567557
$ws.Range("Q7").Value = 1.348258429128
$ws.Range("R7").Value = 12.134325862152
$ws.Range("S7").Value = 0.2530237177323379
$ws.Range("T7").Value = 0.2530237177323379
# Row 8
$ws.Range("G8").Value = 2.260080333333333
$ws.Range("H8").Value = 6.780241
$ws.Range("I8").Value = 0.08199921939034102
$ws.Range("J8").Value = 0.08199921939034102
$ws.Range("M8").Value = 0.06166766666666667
$ws.Range("N8").Value = 0.185003
$ws.Range("O8").Value = 0.3189772891852935
$ws.Range("P8").Value = 0.3189772891852935
$ws.Range("Q8").Value = 0.1393738806358889
$ws.Range("R8").Value = 1.254364925723
$ws.Range("S8").Value = 0.02615588871644113
$ws.Range("T8").Value = 0.02615588871644113
# Row 9
$ws.Range("G9").Value = 2.260080333333333
$ws.Range("H9").Value = 6.780241
$ws.Range("I9").Value = 0.08199921939034102
$ws.Range("J9").Value = 0.08199921939034102
$ws.Range("O9").Value = 0.4045463009579509
$ws.Range("P9").Value = 0.4045463009579509
$ws.Range("Q9").Value = 0.1767623895902222
$ws.Range("R9").Value = 1.590861506312
$ws.Range("S9").Value = 0.03317248088580194
$ws.Range("T9").Value = 0.03317248088580194
# Row 10
$ws.Range("G10").Value = 2.260080333333333
$ws.Range("H10").Value = 6.780241
$ws.Range("I10").Value = 0.08199921939034102
$ws.Range("J10").Value = 0.08199921939034102
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.05345100000000001
$ws.Range("N10").Value = 0.160353
$ws.Range("O10").Value = 0.2764764098567557
$ws.Range("P10").Value = 0.2764764098567557
$ws.Range("Q10").Value = 0.120803553897
$ws.Range("R10").Value = 1.087231985073
$ws.Range("S10").Value = 0.02267084978809795
$ws.Range("T10").Value = 0.02267084978809795
